$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (G=38956)
$ws.Range("H17").Value = 672.4091
$ws.Range("J17").Value = 656.8095
$ws.Range("L17").Value = 1970.4285
$ws.Range("N17").Value = -2306.4285
# Row 33 (G=5512)
$ws.Range("H33").Value = 318.5
$ws.Range("I33").Value = 211.73334
$ws.Range("J33").Value = 496.44446
$ws.Range("K33").Value = 211.73334
$ws.Range("L33").Value = 496.44446
$ws.Range("M33").Value = 17.26666
$ws.Range("N33").Value = -954.4444599999999
# Row 88 (G=12608)
$ws.Range("H88").Value = 5330.25
$ws.Range("I88").Value = 1733
$ws.Range("K88").Value = 1733
$ws.Range("M88").Value = -1327
# Row 91 (G=12608)
$ws.Range("H91").Value = 5330.25
$ws.Range("I91").Value = 1733
$ws.Range("K91").Value = 1733
$ws.Range("M91").Value = -329
# Row 98 (G=36237)
$ws.Range("H98").Value = 6127
$ws.Range("I98").Value = 3813.8333
$ws.Range("K98").Value = 3813.8333
$ws.Range("M98").Value = -2315.8333
# Row 116 (G=27778)
$ws.Range("H116").Value = 19800.445
$ws.Range("I116").Value = 25569.54
$ws.Range("K116").Value = 25569.54
$ws.Range("M116").Value = -22127.54
# Row 122 (G=36237)
$ws.Range("H122").Value = 6127
$ws.Range("I122").Value = 3813.8333
$ws.Range("K122").Value = 11441.4999
$ws.Range("M122").Value = -8991.499899999999
# Row 127 (G=36114)
$ws.Range("H127").Value = 2878.6
$ws.Range("I127").Value = 1500
$ws.Range("J127").Value = 3223.25
$ws.Range("K127").Value = 4500
$ws.Range("L127").Value = 9669.75
$ws.Range("M127").Value = 460
$ws.Range("N127").Value = -19589.75
# Row 131 (G=36108)
$ws.Range("H131").Value = 4281.2383
$ws.Range("I131").Value = 3291.1
$ws.Range("J131").Value = 5181.364
$ws.Range("K131").Value = 9873.299999999999
$ws.Range("L131").Value = 15544.092
$ws.Range("M131").Value = -4833.299999999999
$ws.Range("N131").Value = -25624.092

$ws = $wb.Worksheets.Item("ARM")
# Row 23 (G=2236)
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 20000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -19741
$ws.Range("N23").ClearContents()
# Row 45 (G=27714)
$ws.Range("H45").Value = 3604.2104
$ws.Range("I45").Value = 2051
$ws.Range("K45").Value = 2051
$ws.Range("M45").Value = -1674
# Row 74 (G=44000)
$ws.Range("H74").Value = 321972.4
$ws.Range("I74").Value = 334304.22
$ws.Range("J74").Value = 100000
$ws.Range("K74").Value = 334304.22
$ws.Range("L74").Value = 100000
$ws.Range("M74").Value = -333430.22
$ws.Range("N74").Value = -101748
# Row 77 (G=44000)
$ws.Range("H77").Value = 321972.4
$ws.Range("I77").Value = 334304.22
$ws.Range("J77").Value = 100000
$ws.Range("K77").Value = 1671521.1
$ws.Range("L77").Value = 500000
$ws.Range("M77").Value = -1667153.1
$ws.Range("N77").Value = -508736
# Row 97 (G=19941)
$ws.Range("H97").Value = 1166.4
$ws.Range("I97").Value = 940.3
$ws.Range("K97").Value = 940.3
$ws.Range("M97").Value = -444.3
# Row 102 (G=19945)
$ws.Range("H102").Value = 3077.7273
$ws.Range("I102").Value = 2993.7778
$ws.Range("K102").Value = 2993.7778
$ws.Range("M102").Value = -1371.7778

$ws = $wb.Worksheets.Item("BSM")
# Row 80 (G=13747)
$ws.Range("H80").Value = 744.7857
$ws.Range("I80").Value = 882.6667
$ws.Range("K80").Value = 882.6667
$ws.Range("M80").Value = 115.3333
# Row 83 (G=13747)
$ws.Range("H83").Value = 744.7857
$ws.Range("I83").Value = 882.6667
$ws.Range("K83").Value = 4413.3335
$ws.Range("M83").Value = 578.6665000000003
# Row 99 (G=19943)
$ws.Range("H99").Value = 965
$ws.Range("I99").Value = 965
$ws.Range("K99").Value = 965
$ws.Range("M99").Value = 533
# Row 105 (G=19947)
$ws.Range("H105").Value = 4129.1665
$ws.Range("I105").Value = 3938.25
$ws.Range("K105").Value = 3938.25
$ws.Range("M105").Value = -2191.25

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (G=44021)
$ws.Range("H58").Value = 1273.0344
$ws.Range("I58").Value = 1107.1818
$ws.Range("J58").Value = 1794.2858
$ws.Range("K58").Value = 1107.1818
$ws.Range("L58").Value = 1794.2858
$ws.Range("M58").Value = -904.1818000000001
$ws.Range("N58").Value = -2200.2858
# Row 99 (G=36198)
$ws.Range("H99").Value = 7090.1
$ws.Range("I99").Value = 6749.5
$ws.Range("J99").Value = 7601
$ws.Range("K99").Value = 6749.5
$ws.Range("L99").Value = 7601
$ws.Range("M99").Value = -5251.5
$ws.Range("N99").Value = -10597
# Row 107 (G=27689)
$ws.Range("H107").Value = 1621.125
$ws.Range("I107").Value = 1494.8334
$ws.Range("K107").Value = 1494.8334
$ws.Range("M107").Value = 425.1666
# Row 126 (G=36198)
$ws.Range("H126").Value = 7090.1
$ws.Range("I126").Value = 6749.5
$ws.Range("J126").Value = 7601
$ws.Range("K126").Value = 20248.5
$ws.Range("L126").Value = 22803
$ws.Range("M126").Value = -17778.5
$ws.Range("N126").Value = -27743
# Row 134 (G=44020)
$ws.Range("H134").Value = 2884.6667
$ws.Range("I134").Value = 2270.8096
$ws.Range("K134").Value = 6812.4288
$ws.Range("M134").Value = -4277.4288
# Row 136 (G=44021)
$ws.Range("H136").Value = 1273.0344
$ws.Range("I136").Value = 1107.1818
$ws.Range("J136").Value = 1794.2858
$ws.Range("K136").Value = 3321.5454
$ws.Range("L136").Value = 5382.857400000001
$ws.Range("M136").Value = -771.5454
$ws.Range("N136").Value = -10482.8574

$ws = $wb.Worksheets.Item("CUL")
# Row 137 (G=44088)
$ws.Range("H137").Value = 3652
$ws.Range("J137").Value = 4016.5
$ws.Range("L137").Value = 12049.5
$ws.Range("N137").Value = -22249.5
# Row 140 (G=44097)
$ws.Range("H140").Value = 1926.1904
$ws.Range("I140").Value = 1914.7059
$ws.Range("J140").Value = 1975
$ws.Range("K140").Value = 5744.1177
$ws.Range("L140").Value = 5925
$ws.Range("M140").Value = -564.1176999999998
$ws.Range("N140").Value = -16285
# Row 141 (G=44076)
$ws.Range("H141").Value = 4626.5884
$ws.Range("I141").Value = 3343.25
$ws.Range("K141").Value = 10029.75
$ws.Range("M141").Value = -4849.75

$ws = $wb.Worksheets.Item("GSM")
# Row 43 (G=4218)
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 113 (G=27710)
$ws.Range("H113").Value = 2679.6538
$ws.Range("I113").Value = 2493.2942
$ws.Range("J113").Value = 3031.6667
$ws.Range("K113").Value = 2493.2942
$ws.Range("L113").Value = 3031.6667
$ws.Range("M113").Value = -323.2941999999998
$ws.Range("N113").Value = -7371.6667
# Row 132 (G=44008)
$ws.Range("H132").Value = 2178.5454
$ws.Range("I132").Value = 1151.5
$ws.Range("K132").Value = 3454.5
$ws.Range("M132").Value = -924.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G=36249)
$ws.Range("H7").Value = 4463.5835
$ws.Range("I7").Value = 4646.3
$ws.Range("K7").Value = 4646.3
$ws.Range("M7").Value = -4534.3
# Row 30 (G=1688)
$ws.Range("H30").Value = 6860.8335
$ws.Range("I30").Value = 755.3333
$ws.Range("J30").Value = 12966.333
$ws.Range("K30").Value = 755.3333
$ws.Range("L30").Value = 12966.333
$ws.Range("M30").Value = -647.3333
$ws.Range("N30").Value = -13182.333
# Row 82 (G=12565)
$ws.Range("H82").Value = 2543.6538
$ws.Range("J82").Value = 3817.3845
$ws.Range("L82").Value = 3817.3845
$ws.Range("N82").Value = -4539.3845
# Row 85 (G=12565)
$ws.Range("H85").Value = 2543.6538
$ws.Range("J85").Value = 3817.3845
$ws.Range("L85").Value = 3817.3845
$ws.Range("N85").Value = -6313.3845
# Row 126 (G=36249)
$ws.Range("H126").Value = 4463.5835
$ws.Range("I126").Value = 4646.3
$ws.Range("K126").Value = 13938.9
$ws.Range("M126").Value = -11468.9

$ws = $wb.Worksheets.Item("WVR")
# Row 12 (G=3316)
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
# Row 43 (G=3831)
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
# Row 62 (G=12589)
$ws.Range("H62").Value = 4001.5
$ws.Range("J62").Value = 4001.5
$ws.Range("L62").Value = 4001.5
$ws.Range("N62").Value = -5249.5
# Row 65 (G=12589)
$ws.Range("H65").Value = 4001.5
$ws.Range("J65").Value = 4001.5
$ws.Range("L65").Value = 20007.5
$ws.Range("N65").Value = -26247.5
# Row 81 (G=12596)
$ws.Range("H81").Value = 6170.353
$ws.Range("I81").Value = 6459.7334
$ws.Range("K81").Value = 12919.4668
$ws.Range("M81").Value = -11858.4668
# Row 84 (G=12596)
$ws.Range("H84").Value = 6170.353
$ws.Range("I84").Value = 6459.7334
$ws.Range("K84").Value = 64597.334
$ws.Range("M84").Value = -59293.334
# Row 126 (G=36210)
$ws.Range("H126").Value = 1252999.8
$ws.Range("I126").Value = 997
$ws.Range("K126").Value = 2991
$ws.Range("M126").Value = -521
# Row 132 (G=44029)
$ws.Range("H132").Value = 15671536
$ws.Range("J132").Value = 4306
$ws.Range("L132").Value = 12918
$ws.Range("N132").Value = -17978
# Row 136 (G=44031)
$ws.Range("H136").Value = 15603.341
$ws.Range("I136").Value = 18456.895
$ws.Range("J136").Value = 3555
$ws.Range("K136").Value = 55370.685
$ws.Range("L136").Value = 10665
$ws.Range("M136").Value = -52820.685
$ws.Range("N136").Value = -15765

Write-Output "applied edits"